# Applies the crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.557.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.14%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.482.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.81%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.479.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "

# Row 10
$ws.Range("E10").Value = "  +2.59%  "

# Row 11
$ws.Range("E11").Value = "  +1.18%  "

# Row 12
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.360"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.70%  "

# Row 13
$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "

# Row 14
$ws.Range("E14").Value = "  +2.80%  "

# Row 15
$ws.Range("E15").Value = "  +3.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.924.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.237.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.471.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.78%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.93%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.34%  "

# Row 24
$ws.Range("E24").Value = "  +0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.66%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "638.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.76%  "

# Row 27
$ws.Range("E27").Value = "  +10.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.52%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.593.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.143"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "

# Row 34
$ws.Range("E34").Value = "  +1.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.88%  "

# Row 36
$ws.Range("E36").Value = "  -0.44%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("E38").Value = "  +1.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.98%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.23%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.97%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.64%  "

# Row 44
$ws.Range("E44").Value = "  +0.12%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.80%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0544"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.16%  "

# Row 49
$ws.Range("E49").Value = "  +2.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0236"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.64%  "

# Row 51
$ws.Range("E51").Value = "  +0.06%  "
